$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.458.01'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.918.25'
$ws.Range('E3').Value = '  +1.19%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.67%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.37'
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('E6').Value = '  +0.62%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4829'
$ws.Range('E7').Value = '  +1.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4079'
$ws.Range('E8').Value = '  +0.57%  '
$ws.Range('E9').Value = '  +2.58%  '
$ws.Range('E10').Value = '  +1.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.43'
$ws.Range('E11').Value = '  +0.42%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.911.75'
$ws.Range('E12').Value = '  +2.15%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.081'
$ws.Range('E13').Value = '  +2.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.242'
$ws.Range('E14').Value = '  +2.70%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.48'
$ws.Range('E15').Value = '  +2.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06813'
$ws.Range('E16').Value = '  +2.03%  '
$ws.Range('E17').Value = '  +0.67%  '
$ws.Range('E18').Value = '  +1.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.71'
$ws.Range('E19').Value = '  +0.43%  '
$ws.Range('E20').Value = '  +0.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '29.479.66'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.658'
$ws.Range('E22').Value = '  +2.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.81'
$ws.Range('E23').Value = '  +1.07%  '
$ws.Range('E24').Value = '  +1.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.126.30'
$ws.Range('E25').Value = '  +0.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.644'
$ws.Range('E26').Value = '  +9.87%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '155.96'
$ws.Range('E27').Value = '  +1.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.09'
$ws.Range('E28').Value = '  +1.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.116'
$ws.Range('E29').Value = '  +1.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.53'
$ws.Range('E30').Value = '  +2.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.023'
$ws.Range('E31').Value = '  +0.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09603'
$ws.Range('E32').Value = '  +1.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.689'
$ws.Range('E33').Value = '  +6.06%  '
$ws.Range('E34').Value = '  +0.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.372'
$ws.Range('E35').Value = '  -0.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02289'
$ws.Range('E36').Value = '  +1.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06118'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.182'
$ws.Range('E38').Value = '  +1.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.092'
$ws.Range('E39').Value = '  +3.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5999'
$ws.Range('E40').Value = '  +2.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '10.82'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1849'
$ws.Range('E42').Value = '  +0.46%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.414'
$ws.Range('E43').Value = '  -0.29%  '
$ws.Range('E44').Value = '  -0.46%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.07602'
$ws.Range('E45').Value = '  -1.47%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.43'
$ws.Range('E46').Value = '  +1.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5584'
$ws.Range('E47').Value = '  +1.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.957'
$ws.Range('E48').Value = '  +2.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '118.60'
$ws.Range('E49').Value = '  +5.08%  '
$ws.Range('E50').Value = '  +4.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.51'
$ws.Range('E51').Value = '  +0.95%  '
